$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (pushes current rows 24-27 down to 25-28),
# then fill it with the new weekly price observation for
# "Bruselas (repollito)" in the "Vega Central Mapocho de Santiago" market.
$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value2 = 9
$ws.Cells.Item(24, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value2 = 44463
$ws.Cells.Item(24, 5).Value2 = 13
$ws.Cells.Item(24, 6).Value2 = 100112035
$ws.Cells.Item(24, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value2 = 25
$ws.Cells.Item(24, 11).Value2 = 24000
$ws.Cells.Item(24, 12).Value2 = 25000
$ws.Cells.Item(24, 13).Value2 = 24480
$ws.Cells.Item(24, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(24, 15).Value = "Hijuelas"
$ws.Cells.Item(24, 16).Value2 = 1632
$ws.Cells.Item(24, 17).Value2 = 15
$ws.Cells.Item(24, 18).Value = "Hortaliza"
